$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C2").Value = 287
$ws.Range("C3").Value = 177123
$ws.Range("C4").Value = 167080
$ws.Range("C8").Value = 64.68000000000001
